$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
$data[0,0] = "ECs"
$data[0,1] = "Serpinc1"
$data[0,2] = "Sdc2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 4.551825333333333
$data[0,7] = 13.655476
$data[0,8] = 0.3901863008207799
$data[0,9] = 0.3901863008207799
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 3.057109
$data[0,13] = 9.171327
$data[0,14] = 0.02694952608666365
$data[0,15] = 0.02694952608666365
$data[0,16] = 13.91542619296133
$data[0,17] = 125.238835736652
$data[0,18] = 0.0105153358926284
$data[0,19] = 0.0105153358926284

$data[1,0] = "ECs"
$data[1,1] = "Serpinc1"
$data[1,2] = "Sdc2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 4.551825333333333
$data[1,7] = 13.655476
$data[1,8] = 0.3901863008207799
$data[1,9] = 0.3901863008207799
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 89.02755999999999
$data[1,13] = 267.08268
$data[1,14] = 0.7848102735793893
$data[1,15] = 0.7848102735793893
$data[1,16] = 405.2379029728533
$data[1,17] = 3647.14112675568
$data[1,18] = 0.3062222174940862
$data[1,19] = 0.3062222174940862

$data[2,0] = "ECs"
$data[2,1] = "Serpinc1"
$data[2,2] = "Sdc2"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 4.551825333333333
$data[2,7] = 13.655476
$data[2,8] = 0.3901863008207799
$data[2,9] = 0.3901863008207799
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.184005
$data[2,13] = 0.5520149999999999
$data[2,14] = 0.0016220709001794
$data[2,15] = 0.0016220709001794
$data[2,16] = 0.8375586204599998
$data[2,17] = 7.538027584139999
$data[2,18] = 0.0006329098442100326
$data[2,19] = 0.0006329098442100325

$data[3,0] = "ECs"
$data[3,1] = "Serpinc1"
$data[3,2] = "Sdc2"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 4.551825333333333
$data[3,7] = 13.655476
$data[3,8] = 0.3901863008207799
$data[3,9] = 0.3901863008207799
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 21.16964733333333
$data[3,13] = 63.508942
$data[3,14] = 0.1866181294337677
$data[3,15] = 0.1866181294337677
$data[3,16] = 96.36053702959912
$data[3,17] = 867.244833266392
$data[3,18] = 0.07281583758985534
$data[3,19] = 0.07281583758985534

$data[4,0] = "FAPs"
$data[4,1] = "Serpinc1"
$data[4,2] = "Sdc2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 3.622039333333333
$data[4,7] = 10.866118
$data[4,8] = 0.3104842618962599
$data[4,9] = 0.3104842618962599
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 3.057109
$data[4,13] = 9.171327
$data[4,14] = 0.02694952608666365
$data[4,15] = 0.02694952608666365
$data[4,16] = 11.07296904428733
$data[4,17] = 99.656721398586
$data[4,18] = 0.008367403715471764
$data[4,19] = 0.008367403715471764

$data[5,0] = "FAPs"
$data[5,1] = "Serpinc1"
$data[5,2] = "Sdc2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 3.622039333333333
$data[5,7] = 10.866118
$data[5,8] = 0.3104842618962599
$data[5,9] = 0.3104842618962599
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 89.02755999999999
$data[5,13] = 267.08268
$data[5,14] = 0.7848102735793893
$data[5,15] = 0.7848102735793893
$data[5,16] = 322.4613240706933
$data[5,17] = 2902.15191663624
$data[5,18] = 0.2436712385208985
$data[5,19] = 0.2436712385208985

$data[6,0] = "FAPs"
$data[6,1] = "Serpinc1"
$data[6,2] = "Sdc2"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 3.622039333333333
$data[6,7] = 10.866118
$data[6,8] = 0.3104842618962599
$data[6,9] = 0.3104842618962599
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.184005
$data[6,13] = 0.5520149999999999
$data[6,14] = 0.0016220709001794
$data[6,15] = 0.0016220709001794
$data[6,16] = 0.6664733475299999
$data[6,17] = 5.998260127769999
$data[6,18] = 0.0005036274861856028
$data[6,19] = 0.0005036274861856027

$data[7,0] = "FAPs"
$data[7,1] = "Serpinc1"
$data[7,2] = "Sdc2"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 3.622039333333333
$data[7,7] = 10.866118
$data[7,8] = 0.3104842618962599
$data[7,9] = 0.3104842618962599
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 21.16964733333333
$data[7,13] = 63.508942
$data[7,14] = 0.1866181294337677
$data[7,15] = 0.1866181294337677
$data[7,16] = 76.67729531412844
$data[7,17] = 690.0956578271559
$data[7,18] = 0.05794199217370407
$data[7,19] = 0.05794199217370406

$data[8,0] = "M2"
$data[8,1] = "Serpinc1"
$data[8,2] = "Sdc2"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.291101666666667
$data[8,7] = 3.873305
$data[8,8] = 0.1106743221474397
$data[8,9] = 0.1106743221474397
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 3.057109
$data[8,13] = 9.171327
$data[8,14] = 0.02694952608666365
$data[8,15] = 0.02694952608666365
$data[8,16] = 3.947038525081667
$data[8,17] = 35.523346725735
$data[8,18] = 0.002982620531836243
$data[8,19] = 0.002982620531836242

$data[9,0] = "M2"
$data[9,1] = "Serpinc1"
$data[9,2] = "Sdc2"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.291101666666667
$data[9,7] = 3.873305
$data[9,8] = 0.1106743221474397
$data[9,9] = 0.1106743221474397
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 89.02755999999999
$data[9,13] = 267.08268
$data[9,14] = 0.7848102735793893
$data[9,15] = 0.7848102735793893
$data[9,16] = 114.9436310952667
$data[9,17] = 1034.4926798574
$data[9,18] = 0.0868583450427456
$data[9,19] = 0.08685834504274559

$data[10,0] = "M2"
$data[10,1] = "Serpinc1"
$data[10,2] = "Sdc2"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.291101666666667
$data[10,7] = 3.873305
$data[10,8] = 0.1106743221474397
$data[10,9] = 0.1106743221474397
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.184005
$data[10,13] = 0.5520149999999999
$data[10,14] = 0.0016220709001794
$data[10,15] = 0.0016220709001794
$data[10,16] = 0.237569162175
$data[10,17] = 2.138122459575
$data[10,18] = 0.0001795215973524424
$data[10,19] = 0.0001795215973524423

$data[11,0] = "M2"
$data[11,1] = "Serpinc1"
$data[11,2] = "Sdc2"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.291101666666667
$data[11,7] = 3.873305
$data[11,8] = 0.1106743221474397
$data[11,9] = 0.1106743221474397
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 21.16964733333333
$data[11,13] = 63.508942
$data[11,14] = 0.1866181294337677
$data[11,15] = 0.1866181294337677
$data[11,16] = 27.33216695481223
$data[11,17] = 245.98950259331
$data[11,18] = 0.02065383497550541
$data[11,19] = 0.0206538349755054

$data[12,0] = "sCs"
$data[12,1] = "Serpinc1"
$data[12,2] = "Sdc2"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 2.200808
$data[12,7] = 6.602424000000001
$data[12,8] = 0.1886551151355205
$data[12,9] = 0.1886551151355205
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 3.057109
$data[12,13] = 9.171327
$data[12,14] = 0.02694952608666365
$data[12,15] = 0.02694952608666365
$data[12,16] = 6.728109944072001
$data[12,17] = 60.55298949664801
$data[12,18] = 0.005084165946727245
$data[12,19] = 0.005084165946727245

$data[13,0] = "sCs"
$data[13,1] = "Serpinc1"
$data[13,2] = "Sdc2"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 2.200808
$data[13,7] = 6.602424000000001
$data[13,8] = 0.1886551151355205
$data[13,9] = 0.1886551151355205
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 89.02755999999999
$data[13,13] = 267.08268
$data[13,14] = 0.7848102735793893
$data[13,15] = 0.7848102735793893
$data[13,16] = 195.93256626848
$data[13,17] = 1763.39309641632
$data[13,18] = 0.148058472521659
$data[13,19] = 0.148058472521659

$data[14,0] = "sCs"
$data[14,1] = "Serpinc1"
$data[14,2] = "Sdc2"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 2.200808
$data[14,7] = 6.602424000000001
$data[14,8] = 0.1886551151355205
$data[14,9] = 0.1886551151355205
$data[14,10] = 2
$data[14,11] = 0.6666666666666666
$data[14,12] = 0.184005
$data[14,13] = 0.5520149999999999
$data[14,14] = 0.0016220709001794
$data[14,15] = 0.0016220709001794
$data[14,16] = 0.40495967604
$data[14,17] = 3.64463708436
$data[14,18] = 0.0003060119724313221
$data[14,19] = 0.000306011972431322

$data[15,0] = "sCs"
$data[15,1] = "Serpinc1"
$data[15,2] = "Sdc2"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 2.200808
$data[15,7] = 6.602424000000001
$data[15,8] = 0.1886551151355205
$data[15,9] = 0.1886551151355205
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 21.16964733333333
$data[15,13] = 63.508942
$data[15,14] = 0.1866181294337677
$data[15,15] = 0.1866181294337677
$data[15,16] = 46.59032920837868
$data[15,17] = 419.3129628754081
$data[15,18] = 0.03520646469470293
$data[15,19] = 0.03520646469470292

$ws.Range("A2:T17").Value = $data